$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.6685393258426966
$summary.Range("C2").Value = 0.6056338028169014
$summary.Range("D2").Value = 0.9662921348314607
$summary.Range("E2").Value = 0.7445887445887446
$summary.Range("F2").Value = 0.8634538152610441
$summary.Range("G2").Value = 0.9446556822982678
$summary.Range("H2").Value = 0.798103494227721
$summary.Range("I2").Value = 516
$summary.Range("J2").Value = 336
$summary.Range("K2").Value = 198
$summary.Range("L2").Value = 18

# --- Sheet: Classification Report ---
$report = $wb.Worksheets.Item("Classification Report")
$report.Range("B2").Value = 0.9166666666666666
$report.Range("C2").Value = 0.3707865168539326
$report.Range("D2").Value = 0.528

$report.Range("B3").Value = 0.6056338028169014
$report.Range("C3").Value = 0.9662921348314607
$report.Range("D3").Value = 0.7445887445887446

$report.Range("B4").Value = 0.6685393258426966
$report.Range("C4").Value = 0.6685393258426966
$report.Range("D4").Value = 0.6685393258426966
$report.Range("E4").Value = 0.6685393258426966

$report.Range("B5").Value = 0.761150234741784
$report.Range("C5").Value = 0.6685393258426966
$report.Range("D5").Value = 0.6362943722943724

$report.Range("B6").Value = 0.761150234741784
$report.Range("C6").Value = 0.6685393258426966
$report.Range("D6").Value = 0.6362943722943722

# --- Sheet: Confusion Matrix ---
$confusion = $wb.Worksheets.Item("Confusion Matrix")
$confusion.Range("B2").Value = 198
$confusion.Range("C2").Value = 336
$confusion.Range("B3").Value = 18
$confusion.Range("C3").Value = 516
